$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2019
$ws.Range("J17").Value = 2019
$ws.Range("L17").Value = 6057
$ws.Range("N17").Value = -6393

$ws.Range("H51").Value = 9446.083000000001
$ws.Range("J51").Value = 8690
$ws.Range("L51").Value = 8690
$ws.Range("N51").Value = -9658

$ws.Range("H62").Value = 2610.0588
$ws.Range("I62").Value = 2610.0588
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2610.0588
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -1986.0588
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 2610.0588
$ws.Range("I65").Value = 2610.0588
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 13050.294
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -9930.293999999998
$ws.Range("N65").ClearContents()

$ws.Range("H132").Value = 8195.35
$ws.Range("I132").Value = 8195.35
$ws.Range("K132").Value = 24586.05
$ws.Range("M132").Value = -22056.05

$ws.Range("H138").Value = 323701.28
$ws.Range("I138").Value = 3677
$ws.Range("J138").Value = 483713.44
$ws.Range("K138").Value = 11031
$ws.Range("L138").Value = 1451140.32
$ws.Range("M138").Value = -5891
$ws.Range("N138").Value = -1461420.32

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4450.8335
$ws.Range("J61").Value = 7007.5
$ws.Range("L61").Value = 7007.5
$ws.Range("N61").Value = -7431.5

$ws.Range("H74").Value = 193166.52
$ws.Range("I74").Value = 278930.9
$ws.Range("J74").Value = 2579
$ws.Range("K74").Value = 278930.9
$ws.Range("L74").Value = 2579
$ws.Range("M74").Value = -278056.9
$ws.Range("N74").Value = -4327

$ws.Range("H77").Value = 193166.52
$ws.Range("I77").Value = 278930.9
$ws.Range("J77").Value = 2579
$ws.Range("K77").Value = 1394654.5
$ws.Range("L77").Value = 12895
$ws.Range("M77").Value = -1390286.5
$ws.Range("N77").Value = -21631

$ws.Range("H122").Value = 4051.2104
$ws.Range("I122").Value = 3763.0938
$ws.Range("J122").Value = 5587.8335
$ws.Range("K122").Value = 11289.2814
$ws.Range("L122").Value = 16763.5005
$ws.Range("M122").Value = -8839.2814
$ws.Range("N122").Value = -21663.5005

$ws.Range("H132").Value = 3137.5
$ws.Range("I132").Value = 1942.65
$ws.Range("J132").Value = 6124.625
$ws.Range("K132").Value = 5827.950000000001
$ws.Range("L132").Value = 18373.875
$ws.Range("M132").Value = -3297.950000000001
$ws.Range("N132").Value = -23433.875

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws.Range("H136").Value = 4450.8335
$ws.Range("J136").Value = 7007.5
$ws.Range("L136").Value = 21022.5
$ws.Range("N136").Value = -26122.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2971.1428
$ws.Range("I134").Value = 2736.2666
$ws.Range("J134").Value = 3558.3333
$ws.Range("K134").Value = 8208.799800000001
$ws.Range("L134").Value = 10674.9999
$ws.Range("M134").Value = -5673.799800000001
$ws.Range("N134").Value = -15744.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6808.9414
$ws.Range("I31").Value = 5615.727
$ws.Range("J31").Value = 8996.5
$ws.Range("K31").Value = 5615.727
$ws.Range("L31").Value = 8996.5
$ws.Range("M31").Value = -5320.727
$ws.Range("N31").Value = -9586.5

$ws.Range("H34").Value = 6808.9414
$ws.Range("I34").Value = 5615.727
$ws.Range("J34").Value = 8996.5
$ws.Range("K34").Value = 5615.727
$ws.Range("L34").Value = 8996.5
$ws.Range("M34").Value = -5413.727
$ws.Range("N34").Value = -9400.5

$ws.Range("H132").Value = 4666.263
$ws.Range("J132").Value = 5249.25
$ws.Range("L132").Value = 15747.75
$ws.Range("N132").Value = -20807.75

$ws.Range("H134").Value = 4795.385
$ws.Range("I134").Value = 4795.385
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 14386.155
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -11851.155
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 1093.3889
$ws.Range("J44").Value = 1112.0667
$ws.Range("L44").Value = 3336.2001
$ws.Range("N44").Value = -4132.2001

$ws.Range("H113").Value = 5755.5415
$ws.Range("I113").Value = 448.42856
$ws.Range("K113").Value = 1345.28568
$ws.Range("M113").Value = 824.71432

$ws.Range("H122").Value = 1621.1875
$ws.Range("J122").Value = 1940.8182
$ws.Range("L122").Value = 17467.3638
$ws.Range("N122").Value = -22367.3638

$ws.Range("H124").Value = 3416
$ws.Range("I124").Value = 3416
$ws.Range("K124").Value = 10248
$ws.Range("M124").Value = -5338

$ws.Range("H125").Value = 2000
$ws.Range("I125").Value = 2000
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 6000
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -1080
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 158728.16
$ws.Range("I70").Value = 204549.6
$ws.Range("K70").Value = 204549.6
$ws.Range("M70").Value = -204279.6

$ws.Range("H73").Value = 158728.16
$ws.Range("I73").Value = 204549.6
$ws.Range("K73").Value = 204549.6
$ws.Range("M73").Value = -203613.6

$ws.Range("H102").Value = 1555.9474
$ws.Range("I102").Value = 797.1667
$ws.Range("J102").Value = 2856.7144
$ws.Range("K102").Value = 797.1667
$ws.Range("L102").Value = 2856.7144
$ws.Range("M102").Value = 824.8333
$ws.Range("N102").Value = -6100.7144

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3373.75
$ws.Range("I122").Value = 3355.7144
$ws.Range("K122").Value = 10067.1432
$ws.Range("M122").Value = -7617.143199999999

$ws.Range("H132").Value = 8159.5
$ws.Range("I132").Value = 4487
$ws.Range("J132").Value = 13301
$ws.Range("K132").Value = 13461
$ws.Range("L132").Value = 39903
$ws.Range("M132").Value = -10931
$ws.Range("N132").Value = -44963

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 53918.43
$ws.Range("J46").Value = 53918.43
$ws.Range("L46").Value = 53918.43
$ws.Range("N46").Value = -54380.43

$ws.Range("H134").Value = 53918.43
$ws.Range("J134").Value = 53918.43
$ws.Range("L134").Value = 161755.29
$ws.Range("N134").Value = -166825.29
